$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix data values that changed
$ws.Range("D2").Value = 28
$ws.Range("E6").Value = 3

# Average Score formulas (F2:F10) - enter first two individually, then fill the rest
# to mirror the authoring process (matches shared-formula grouping in target file)
$ws.Range("F2").Formula = "=AVERAGE(C2:E2)"
$ws.Range("F3").Formula = "=AVERAGE(C3:E3)"
$ws.Range("F4:F10").Formula = "=AVERAGE(C4:E4)"

# Final Score formulas (H2:H10) - enter first individually, then fill the rest
$ws.Range("H2").Formula = "=SUM(F2,G2)"
$ws.Range("H3:H10").Formula = "=SUM(F3,G3)"

# Summary statistics
$ws.Range("B15").Formula = "=AVERAGE(H2:H10)"
$ws.Range("B16").Formula = "=MEDIAN(H2:H10)"
$ws.Range("B17").Formula = "=MAX(H2:H10)"
$ws.Range("B18").Formula = "=MIN(H2:H10)"
$ws.Range("B19").Formula = "=STDEV(H1:H10)"

$wb.Save()
